$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert 11 new rows starting at row 8 (pushes the existing row 8..21 block down to 19..32),
# mirroring the blank-row-separated groups already used on this sheet.
$null = $ws.Range("A8:A18").EntireRow.Insert()

# Populate the newly inserted rows with the new log-message Name/Value pairs.
$ws.Range("A9").Value = "logMessageExtractDataJob"
$ws.Range("B9").Value = "Starting the data extraction of jobs download folder…"

$ws.Range("A10").Value = "logMessageExtractDataCV"
$ws.Range("B10").Value = "Starting the data extraction of CVs download folder…"

$ws.Range("A12").Value = "logMessageDocumentUnderstanding"
$ws.Range("B12").Value = "Starting Document Understanding: "

$ws.Range("A14").Value = "logMessageMoveFile"
$ws.Range("B14").Value = "Moving the file to the archive: "

$ws.Range("A16").Value = "logGenerateExcelFile"
$ws.Range("B16").Value = "Generating the Excel File: "

$ws.Range("A18").Value = "logSendingEmail"
$ws.Range("B18").Value = "Sending the email…"

# Column A got a touch wider to fit the new, longer constant names.
$ws.Columns.Item(1).ColumnWidth = 30.4973958333333

# Match the saved selection/active cell from the authored workbook.
$null = $ws.Range("B18").Select()
